$d = $word.ActiveDocument

# Table 3 = "ХРОНОЛОГИЯ НА СЪБИТИЯТА" (chronology table, 3 columns: N / Компания / Дата и час)
# Rows 5-8 (1-based, i.e. entries N=3..N=6) have empty "Компания" / "Дата" cells that
# need to be filled in with the new chronology entries.

function Set-CellText($rowIndex, $colIndex, $text) {
    $t = $d.Tables.Item(3)
    $cell = $t.Cell($rowIndex, $colIndex)
    $cell.Range.Text = $text
    $t = $d.Tables.Item(3)
    $cell = $t.Cell($rowIndex, $colIndex)
    $cell.Range.Font.Name = "Calibri"
}

Set-CellText 5 2 "KILL THE WORLD"
Set-CellText 5 3 "21.10.2021 / 14:39"

Set-CellText 6 2 "KILL THE WORLD"
Set-CellText 6 3 "23.10.2021 / 21:49"

Set-CellText 7 2 "CHINA INC."
Set-CellText 7 3 "21.10.2021 / 14:48"

Set-CellText 8 2 "CHINA INC."
Set-CellText 8 3 "23.10.2021 / 21:56"

Write-Output "chronology rows filled"
